$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 21,20

$arr[0,0] = "Sending cluster"
$arr[0,1] = "Ligand symbol"
$arr[0,2] = "Receptor symbol"
$arr[0,3] = "Target cluster"
$arr[0,4] = "Ligand-expressing cells"
$arr[0,5] = "Ligand detection rate"
$arr[0,6] = "Ligand average expression value"
$arr[0,7] = "Ligand total expression value"
$arr[0,8] = "Ligand derived specificity of average expression value"
$arr[0,9] = "Ligand derived specificity of total expression value"
$arr[0,10] = "Receptor-expressing cells"
$arr[0,11] = "Receptor detection rate"
$arr[0,12] = "Receptor average expression value"
$arr[0,13] = "Receptor total expression value"
$arr[0,14] = "Receptor derived specificity of average expression value"
$arr[0,15] = "Receptor derived specificity of total expression value"
$arr[0,16] = "Edge average expression weight"
$arr[0,17] = "Edge total expression weight"
$arr[0,18] = "Edge average expression derived specificity"
$arr[0,19] = "Edge total expression derived specificity"

$arr[1,0] = "ECs"
$arr[1,1] = "Inhbb"
$arr[1,2] = "Acvr1b"
$arr[1,3] = "ECs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 4.308396
$arr[1,7] = 12.925188
$arr[1,8] = 0.8865596696737007
$arr[1,9] = 0.8865596696737006
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 3.479406
$arr[1,13] = 10.438218
$arr[1,14] = 0.1624220085658938
$arr[1,15] = 0.1624220085658938
$arr[1,16] = 14.990658892776
$arr[1,17] = 134.915930034984
$arr[1,18] = 0.1439968022619177
$arr[1,19] = 0.1439968022619177

$arr[2,0] = "ECs"
$arr[2,1] = "Inhbb"
$arr[2,2] = "Acvr1b"
$arr[2,3] = "FAPs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 4.308396
$arr[2,7] = 12.925188
$arr[2,8] = 0.8865596696737007
$arr[2,9] = 0.8865596696737006
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 7.930816666666668
$arr[2,13] = 23.79245
$arr[2,14] = 0.3702181270503834
$arr[2,15] = 0.3702181270503834
$arr[2,16] = 34.1690988034
$arr[2,17] = 307.5218892306
$arr[2,18] = 0.328220460425004
$arr[2,19] = 0.328220460425004

$arr[3,0] = "ECs"
$arr[3,1] = "Inhbb"
$arr[3,2] = "Acvr1b"
$arr[3,3] = "Inflammatory-Mac"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 4.308396
$arr[3,7] = 12.925188
$arr[3,8] = 0.8865596696737007
$arr[3,9] = 0.8865596696737006
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 4.767182666666667
$arr[3,13] = 14.301548
$arr[3,14] = 0.2225366582458366
$arr[3,15] = 0.2225366582458366
$arr[3,16] = 20.538910732336
$arr[3,17] = 184.850196591024
$arr[3,18] = 0.1972920262247181
$arr[3,19] = 0.1972920262247181

$arr[4,0] = "ECs"
$arr[4,1] = "Inhbb"
$arr[4,2] = "Acvr1b"
$arr[4,3] = "MuSCs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 4.308396
$arr[4,7] = 12.925188
$arr[4,8] = 0.8865596696737007
$arr[4,9] = 0.8865596696737006
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 1.349
$arr[4,13] = 4.047
$arr[4,14] = 0.06297261358846615
$arr[4,15] = 0.06297261358846615
$arr[4,16] = 5.812026204
$arr[4,17] = 52.30823583599999
$arr[4,18] = 0.05582897950148014
$arr[4,19] = 0.05582897950148014

$arr[5,0] = "ECs"
$arr[5,1] = "Inhbb"
$arr[5,2] = "Acvr1b"
$arr[5,3] = "Resolving-Mac"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 4.308396
$arr[5,7] = 12.925188
$arr[5,8] = 0.8865596696737007
$arr[5,9] = 0.8865596696737006
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 3.895605333333334
$arr[5,13] = 11.686816
$arr[5,14] = 0.1818505925494202
$arr[5,15] = 0.1818505925494202
$arr[5,16] = 16.783810435712
$arr[5,17] = 151.054293921408
$arr[5,18] = 0.1612214012605807
$arr[5,19] = 0.1612214012605807

$arr[6,0] = "FAPs"
$arr[6,1] = "Inhbb"
$arr[6,2] = "Acvr1b"
$arr[6,3] = "ECs"
$arr[6,4] = 1
$arr[6,5] = 0.3333333333333333
$arr[6,6] = 0.062595
$arr[6,7] = 0.187785
$arr[6,8] = 0.01288047861042144
$arr[6,9] = 0.01288047861042144
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 3.479406
$arr[6,13] = 10.438218
$arr[6,14] = 0.1624220085658938
$arr[6,15] = 0.1624220085658938
$arr[6,16] = 0.21779341857
$arr[6,17] = 1.96014076713
$arr[6,18] = 0.002092073207194683
$arr[6,19] = 0.002092073207194682

$arr[7,0] = "FAPs"
$arr[7,1] = "Inhbb"
$arr[7,2] = "Acvr1b"
$arr[7,3] = "FAPs"
$arr[7,4] = 1
$arr[7,5] = 0.3333333333333333
$arr[7,6] = 0.062595
$arr[7,7] = 0.187785
$arr[7,8] = 0.01288047861042144
$arr[7,9] = 0.01288047861042144
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 7.930816666666668
$arr[7,13] = 23.79245
$arr[7,14] = 0.3702181270503834
$arr[7,15] = 0.3702181270503834
$arr[7,16] = 0.49642946925
$arr[7,17] = 4.46786522325
$arr[7,18] = 0.004768586666662751
$arr[7,19] = 0.00476858666666275

$arr[8,0] = "FAPs"
$arr[8,1] = "Inhbb"
$arr[8,2] = "Acvr1b"
$arr[8,3] = "Inflammatory-Mac"
$arr[8,4] = 1
$arr[8,5] = 0.3333333333333333
$arr[8,6] = 0.062595
$arr[8,7] = 0.187785
$arr[8,8] = 0.01288047861042144
$arr[8,9] = 0.01288047861042144
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 4.767182666666667
$arr[8,13] = 14.301548
$arr[8,14] = 0.2225366582458366
$arr[8,15] = 0.2225366582458366
$arr[8,16] = 0.29840179902
$arr[8,17] = 2.68561619118
$arr[8,18] = 0.002866378666570165
$arr[8,19] = 0.002866378666570165

$arr[9,0] = "FAPs"
$arr[9,1] = "Inhbb"
$arr[9,2] = "Acvr1b"
$arr[9,3] = "MuSCs"
$arr[9,4] = 1
$arr[9,5] = 0.3333333333333333
$arr[9,6] = 0.062595
$arr[9,7] = 0.187785
$arr[9,8] = 0.01288047861042144
$arr[9,9] = 0.01288047861042144
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 1.349
$arr[9,13] = 4.047
$arr[9,14] = 0.06297261358846615
$arr[9,15] = 0.06297261358846615
$arr[9,16] = 0.08444065499999999
$arr[9,17] = 0.759965895
$arr[9,18] = 0.0008111174023685728
$arr[9,19] = 0.0008111174023685727

$arr[10,0] = "FAPs"
$arr[10,1] = "Inhbb"
$arr[10,2] = "Acvr1b"
$arr[10,3] = "Resolving-Mac"
$arr[10,4] = 1
$arr[10,5] = 0.3333333333333333
$arr[10,6] = 0.062595
$arr[10,7] = 0.187785
$arr[10,8] = 0.01288047861042144
$arr[10,9] = 0.01288047861042144
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 3.895605333333334
$arr[10,13] = 11.686816
$arr[10,14] = 0.1818505925494202
$arr[10,15] = 0.1818505925494202
$arr[10,16] = 0.24384541584
$arr[10,17] = 2.19460874256
$arr[10,18] = 0.002342322667625272
$arr[10,19] = 0.002342322667625271

$arr[11,0] = "Inflammatory-Mac"
$arr[11,1] = "Inhbb"
$arr[11,2] = "Acvr1b"
$arr[11,3] = "ECs"
$arr[11,4] = 1
$arr[11,5] = 0.3333333333333333
$arr[11,6] = 0.05725033333333333
$arr[11,7] = 0.171751
$arr[11,8] = 0.01178068046871951
$arr[11,9] = 0.01178068046871951
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 3.479406
$arr[11,13] = 10.438218
$arr[11,14] = 0.1624220085658938
$arr[11,15] = 0.1624220085658938
$arr[11,16] = 0.199197153302
$arr[11,17] = 1.792774379718
$arr[11,18] = 0.001913441784002417
$arr[11,19] = 0.001913441784002417

$arr[12,0] = "Inflammatory-Mac"
$arr[12,1] = "Inhbb"
$arr[12,2] = "Acvr1b"
$arr[12,3] = "FAPs"
$arr[12,4] = 1
$arr[12,5] = 0.3333333333333333
$arr[12,6] = 0.05725033333333333
$arr[12,7] = 0.171751
$arr[12,8] = 0.01178068046871951
$arr[12,9] = 0.01178068046871951
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 7.930816666666668
$arr[12,13] = 23.79245
$arr[12,14] = 0.3702181270503834
$arr[12,15] = 0.3702181270503834
$arr[12,16] = 0.4540418977722223
$arr[12,17] = 4.08637707995
$arr[12,18] = 0.004361421458508368
$arr[12,19] = 0.004361421458508368

$arr[13,0] = "Inflammatory-Mac"
$arr[13,1] = "Inhbb"
$arr[13,2] = "Acvr1b"
$arr[13,3] = "Inflammatory-Mac"
$arr[13,4] = 1
$arr[13,5] = 0.3333333333333333
$arr[13,6] = 0.05725033333333333
$arr[13,7] = 0.171751
$arr[13,8] = 0.01178068046871951
$arr[13,9] = 0.01178068046871951
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 4.767182666666667
$arr[13,13] = 14.301548
$arr[13,14] = 0.2225366582458366
$arr[13,15] = 0.2225366582458366
$arr[13,16] = 0.2729227967275555
$arr[13,17] = 2.456305170548
$arr[13,18] = 0.002621633263370836
$arr[13,19] = 0.002621633263370835

$arr[14,0] = "Inflammatory-Mac"
$arr[14,1] = "Inhbb"
$arr[14,2] = "Acvr1b"
$arr[14,3] = "MuSCs"
$arr[14,4] = 1
$arr[14,5] = 0.3333333333333333
$arr[14,6] = 0.05725033333333333
$arr[14,7] = 0.171751
$arr[14,8] = 0.01178068046871951
$arr[14,9] = 0.01178068046871951
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 1.349
$arr[14,13] = 4.047
$arr[14,14] = 0.06297261358846615
$arr[14,15] = 0.06297261358846615
$arr[14,16] = 0.07723069966666665
$arr[14,17] = 0.6950762969999998
$arr[14,18] = 0.0007418602389658638
$arr[14,19] = 0.0007418602389658637

$arr[15,0] = "Inflammatory-Mac"
$arr[15,1] = "Inhbb"
$arr[15,2] = "Acvr1b"
$arr[15,3] = "Resolving-Mac"
$arr[15,4] = 1
$arr[15,5] = 0.3333333333333333
$arr[15,6] = 0.05725033333333333
$arr[15,7] = 0.171751
$arr[15,8] = 0.01178068046871951
$arr[15,9] = 0.01178068046871951
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 3.895605333333334
$arr[15,13] = 11.686816
$arr[15,14] = 0.1818505925494202
$arr[15,15] = 0.1818505925494202
$arr[15,16] = 0.2230247038684444
$arr[15,17] = 2.007222334816
$arr[15,18] = 0.002142323723872024
$arr[15,19] = 0.002142323723872024

$arr[16,0] = "MuSCs"
$arr[16,1] = "Inhbb"
$arr[16,2] = "Acvr1b"
$arr[16,3] = "ECs"
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 0.4314383333333334
$arr[16,7] = 1.294315
$arr[16,8] = 0.08877917124715834
$arr[16,9] = 0.08877917124715833
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 3.479406
$arr[16,13] = 10.438218
$arr[16,14] = 0.1624220085658938
$arr[16,15] = 0.1624220085658938
$arr[16,16] = 1.50114912563
$arr[16,17] = 13.51034213067
$arr[16,18] = 0.0144196913127789
$arr[16,19] = 0.0144196913127789

$arr[17,0] = "MuSCs"
$arr[17,1] = "Inhbb"
$arr[17,2] = "Acvr1b"
$arr[17,3] = "FAPs"
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 0.4314383333333334
$arr[17,7] = 1.294315
$arr[17,8] = 0.08877917124715834
$arr[17,9] = 0.08877917124715833
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = 7.930816666666668
$arr[17,13] = 23.79245
$arr[17,14] = 0.3702181270503834
$arr[17,15] = 0.3702181270503834
$arr[17,16] = 3.421658324638889
$arr[17,17] = 30.79492492175001
$arr[17,18] = 0.03286765850020821
$arr[17,19] = 0.0328676585002082

$arr[18,0] = "MuSCs"
$arr[18,1] = "Inhbb"
$arr[18,2] = "Acvr1b"
$arr[18,3] = "Inflammatory-Mac"
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 0.4314383333333334
$arr[18,7] = 1.294315
$arr[18,8] = 0.08877917124715834
$arr[18,9] = 0.08877917124715833
$arr[18,10] = 3
$arr[18,11] = 1
$arr[18,12] = 4.767182666666667
$arr[18,13] = 14.301548
$arr[18,14] = 0.2225366582458366
$arr[18,15] = 0.2225366582458366
$arr[18,16] = 2.056745344402223
$arr[18,17] = 18.51070809962
$arr[18,18] = 0.01975662009117748
$arr[18,19] = 0.01975662009117748

$arr[19,0] = "MuSCs"
$arr[19,1] = "Inhbb"
$arr[19,2] = "Acvr1b"
$arr[19,3] = "MuSCs"
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 0.4314383333333334
$arr[19,7] = 1.294315
$arr[19,8] = 0.08877917124715834
$arr[19,9] = 0.08877917124715833
$arr[19,10] = 3
$arr[19,11] = 1
$arr[19,12] = 1.349
$arr[19,13] = 4.047
$arr[19,14] = 0.06297261358846615
$arr[19,15] = 0.06297261358846615
$arr[19,16] = 0.5820103116666667
$arr[19,17] = 5.238092805
$arr[19,18] = 0.005590656445651566
$arr[19,19] = 0.005590656445651565

$arr[20,0] = "MuSCs"
$arr[20,1] = "Inhbb"
$arr[20,2] = "Acvr1b"
$arr[20,3] = "Resolving-Mac"
$arr[20,4] = 3
$arr[20,5] = 1
$arr[20,6] = 0.4314383333333334
$arr[20,7] = 1.294315
$arr[20,8] = 0.08877917124715834
$arr[20,9] = 0.08877917124715833
$arr[20,10] = 3
$arr[20,11] = 1
$arr[20,12] = 3.895605333333334
$arr[20,13] = 11.686816
$arr[20,14] = 0.1818505925494202
$arr[20,15] = 0.1818505925494202
$arr[20,16] = 1.680713472337778
$arr[20,17] = 15.12642125104
$arr[20,18] = 0.01614454489734219
$arr[20,19] = 0.01614454489734219

$ws.Range("A1:T21").Value = $arr
